$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios")

# Insert a new column before column C ("number_of_run"), shifting the
# existing agent_account_min..account columns one place to the right.
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "number_of_run"

# Data for the new column (3 scenario rows).
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 1

# Match the new column's width to the rest of the "best fit" header row
# (pre-compensates for the engine's pixel-rounding so the stored value
# lands on 14 characters).
$ws.Columns("C").ColumnWidth = 13.285714285714286

# Make "scenarios" the active sheet/tab with the selection left on D4,
# as it was after editing the sheet.
$ws.Activate() | Out-Null
$ws.Range("D4").Select() | Out-Null

# Give the scenarios sheet its own (default) page setup, like the other
# sheets in the workbook already have.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "done"
